# Updated CHE_grids model - 2025-08-15 01:17
# Fixes the grid_cell (column AG) assignments on the "solar" worksheet so that
# each SubRES_New_RE_and_Conventional row points at its corrected CHE grid cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$gridCellUpdates = @{
    4  = "CHE_10"
    5  = "CHE_22"
    6  = "CHE_0"
    8  = "CHE_3"
    9  = "CHE_14"
    10 = "CHE_18"
    11 = "CHE_12"
    12 = "CHE_24"
    13 = "CHE_8"
    14 = "CHE_5"
    15 = "CHE_11"
    16 = "CHE_15"
    17 = "CHE_25"
    18 = "CHE_13"
    19 = "CHE_2"
    20 = "CHE_9"
    21 = "CHE_21"
    22 = "CHE_4"
    26 = "CHE_20"
    27 = "CHE_1"
    28 = "CHE_6"
}

foreach ($row in $gridCellUpdates.Keys) {
    $ws.Range("AG" + $row).Value = $gridCellUpdates[$row]
}
